# Cover-slide date update: "June 2024" -> "July 2024"
#
# The cover slide has a byline textbox ("TextBox 1") whose second paragraph
# reads "June 2024" as a single run. The author retyped just the word
# "June " (keeping the trailing "2024"), which PowerPoint represents as the
# paragraph splitting into two runs: the freshly-typed "July " and the
# untouched "2024" tail.

$p  = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# Find the byline shape on the cover slide: the one whose text contains
# "June 2024" (name is "TextBox 1", but search by content to be robust).
$byline = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $candidate = $s1.Shapes.Item($i)
    if ($candidate.HasTextFrame) {
        if ($candidate.TextFrame.TextRange.Text -like "*June 2024*") {
            $byline = $candidate
        }
    }
}

$tr = $byline.TextFrame.TextRange
$fullText = $tr.Text

# Locate "June " (the word plus trailing space) within the shape's flat
# text so we only retype that portion, leaving "2024" untouched - this is
# what preserves the run split PowerPoint itself would produce.
$needle = "June "
$startIndex0 = $fullText.IndexOf($needle)

# TextRange.Characters uses 1-based, inclusive character positions.
$start = $startIndex0 + 1
$length = $needle.Length

$monthRun = $tr.Characters($start, $length)
$monthRun.Text = "July "
